$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume figures (and the row-31/32 swap)
# Each Price/Volume column is stored as text, so force text format before writing
# to avoid Excel auto-converting numeric-looking strings (e.g. "1.001", "4.400").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.243.08'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.36%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.859.73'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.58%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7011'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '237.38'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.08143'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +9.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3023'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.61%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08158'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.862.38'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.152'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7038'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.77%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '88.99'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.258.96'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.766'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007835'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.30'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.79%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '235.23'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.28%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.111.95'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.50%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.404'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '161.35'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.19%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.940'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1442'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.05'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.957'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.433'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.80%  '
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.400'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.22%  '
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.479'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.99%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05182'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.165'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.60%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7053'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.96%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9979'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.676'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01834'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.76%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9218'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.135.35'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +4.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4255'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.880'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '70.08'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.58%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.16'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.762'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.007.18'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.52%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.144'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.930'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.39%  '
